$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the team record columns, matching the style
# of the existing header row (row 1) by copying formats from AB1.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the team's record (Wins/Losses/Ties) for every player row.
$wins = 91
$losses = 71
$ties = 0

for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 29).Value = $wins    # AC
    $ws.Cells.Item($row, 30).Value = $losses  # AD
    $ws.Cells.Item($row, 31).Value = $ties    # AE
}
